# Updates cryptos list cell values (Price column D, Volume(1h) column E)
# per commit "Updated cryptos list on Mon Apr 29 21:36:17 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text can never be misread as a number (dotted/percent/spaced strings) -
# safe to assign directly.
$directValues = [ordered]@{
    "D2" = "62.943.45"
    "E2" = "  -1.28%  "
    "D3" = "3.168.40"
    "E3" = "  -4.14%  "
    "E4" = "  -0.07%  "
    "E5" = "  -2.71%  "
    "E6" = "  -5.18%  "
    "E7" = "  +0.04%  "
    "D8" = "3.162.91"
    "E8" = "  -4.32%  "
    "E9" = "  -0.88%  "
    "E10" = "  -6.62%  "
    "E11" = "  -6.76%  "
    "E12" = "  -3.98%  "
    "E13" = "  -5.07%  "
    "E14" = "  -2.19%  "
    "D15" = "3.693.27"
    "E15" = "  -4.16%  "
    "E16" = "  -1.69%  "
    "D17" = "3.170.22"
    "D18" = "62.948.63"
    "E18" = "  -1.38%  "
    "E19" = "  -5.38%  "
    "E20" = "  -4.70%  "
    "E21" = "  -0.56%  "
    "E22" = "  -6.45%  "
    "E23" = "  -5.55%  "
    "E24" = "  -5.22%  "
    "E25" = "  -3.69%  "
    "E26" = "  -0.07%  "
    "E27" = "  -0.14%  "
    "E28" = "  -4.41%  "
    "E29" = "  -7.64%  "
    "E30" = "  -7.14%  "
    "E31" = "  -6.17%  "
    "E32" = "  -6.09%  "
    "E33" = "  -4.64%  "
    "E34" = "  -7.04%  "
    "E35" = "  -6.56%  "
    "E36" = "  -5.00%  "
    "E37" = "  -2.52%  "
    "D38" = (@("0.0", [string][char]0x2083, "0702") -join "")
    "E38" = "  -6.46%  "
    "E39" = "  -3.65%  "
    "E40" = "  -7.28%  "
    "E41" = "  -3.19%  "
    "E42" = "  -4.59%  "
    "E43" = "  -6.64%  "
    "D44" = "2.812.57"
    "E44" = "  -9.59%  "
    "E45" = "  -5.76%  "
    "E47" = "  -6.14%  "
    "E48" = "  -1.38%  "
    "E49" = "  -6.38%  "
    "E50" = "  -5.40%  "
    "E51" = "  -2.48%  "
}

foreach ($ref in $directValues.Keys) {
    $ws.Range($ref).Value = $directValues[$ref]
}

# Cells whose new text parses as a plain number (e.g. "5.18", "1.00") - assigning those
# straight to .Value would store a Double and lose the original text formatting, so build
# each one as a quoted-text formula first, then collapse the formula down to its literal
# text value via copy / paste-special (values only). That keeps the cell a plain text cell,
# matching how the source data is stored.
$textValues = [ordered]@{
    "D5" = "590.93"
    "D6" = "134.17"
    "D9" = "0.513"
    "D11" = "5.18"
    "D12" = "0.451"
    "D13" = "0.0000235"
    "D14" = "34.18"
    "D19" = "6.50"
    "D20" = "458.17"
    "D21" = "13.96"
    "D22" = "0.693"
    "D23" = "7.55"
    "D24" = "13.22"
    "D25" = "82.17"
    "D27" = "1.00"
    "D29" = "6.70"
    "D30" = "7.60"
    "D32" = "27.08"
    "D34" = "2.34"
    "D35" = "1.03"
    "D36" = "5.76"
    "D37" = "51.08"
    "D39" = "0.0385"
    "D40" = "400.64"
    "D41" = "8.06"
    "D42" = "2.62"
    "D45" = "0.250"
    "D47" = "2.09"
    "D48" = "124.36"
    "D49" = "34.49"
    "D50" = "24.98"
}

foreach ($ref in $textValues.Keys) {
    $cell = $ws.Range($ref)
    $literal = $textValues[$ref]
    $cell.Formula = "=""$literal"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

